$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "AUC on training set" ---
$ws.Range("B13").Value = 0.85573675100067004
$ws.Range("C13").Value = 0.85773064727643067
$ws.Range("D13").Value = 0.85848441513750184
$ws.Range("E13").Value = 0.86119631930388285
$ws.Range("F13").Value = 0.86280632349525188
$ws.Range("G13").Value = 0.86275591869981272
$ws.Range("H13").Value = 0.86331819516928354
$ws.Range("I13").Value = 0.86366631572375652
$ws.Range("J13").Value = 0.86243397304631564
$ws.Range("K13").Value = 0.86154233992522722

# --- Row 15: "Number observations test set" ---
$ws.Range("B15").Value = 416408
$ws.Range("C15").Value = 427741
$ws.Range("D15").Value = 435403
$ws.Range("E15").Value = 444857
$ws.Range("F15").Value = 458164
$ws.Range("G15").Value = 474323
$ws.Range("H15").Value = 492119
$ws.Range("I15").Value = 509898
$ws.Range("J15").Value = 528492
$ws.Range("K15").Value = 546956

# --- Row 16: "Number bankrupt in test set" ---
$ws.Range("B16").Value = 7331
$ws.Range("C16").Value = 7485
$ws.Range("D16").Value = 6898
$ws.Range("E16").Value = 6645
$ws.Range("F16").Value = 6515
$ws.Range("G16").Value = 6560
$ws.Range("H16").Value = 6447
$ws.Range("I16").Value = 6481
$ws.Range("J16").Value = 6574
$ws.Range("K16").Value = 6924

# Reset the selection back to A1 (matches the saved workbook having no stale
# "H21" selection left over from the previous session).
$ws.Range("A1").Select()

Write-Host "Edits applied"
